$d = $word.ActiveDocument

# Find the unique text " ore” " (smart closing quote + trailing space) that
# ends the "Bring me that sweet sweet ore" " paragraph. Keep a reference to
# the very Range object that Find.Execute mutates, since Find collapses /
# repositions that same Range to the match.
$rng = $d.Content
$found = $rng.Find.Execute("ore” ")
if (-not $found) {
    throw "Could not find target text 'ore” ' in document"
}

$targetPara = $rng.Paragraphs(1)

# Re-insert the "_GoBack" bookmark right after the closing quote, i.e. one
# character before the end of the matched range (which currently sits just
# past the trailing space / at the paragraph mark).
$bmPos = $rng.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the following three paragraphs: the blank paragraph, the
# "Multiple different diggers..." paragraph and the "Need to swap..."
# paragraph, restoring the document to end right after the target
# paragraph (which now ends in "...ore”" + bookmark + a single space run).
$emptyPara = $targetPara.Next()
$lastParaToRemove = $emptyPara.Next().Next()

$deleteStart = $targetPara.Range.End
$deleteEnd = $lastParaToRemove.Range.End
$killRange = $d.Range($deleteStart, $deleteEnd)
$killRange.Delete()
